$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.33%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'40.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.51%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.111"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.41%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07622"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.65%"
$ws.Range("E5").Style = "Normal"
$ws.Range("B6").Value = "'GateToken"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'4.250"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.06%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "'FTXToken"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'1.607"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-1.29%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'2.491"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'2.18%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "'MXToken"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.9042"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.63%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.1108"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'7.71%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'WazirX"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.1781"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.96%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'MandalaExchangeToken"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.09183"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.57%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'BitrueCoin"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.04155"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-5.53%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'BitMartToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.1054"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.24%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'BitForexToken"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.001253"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.64%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'TigerCash"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.005679"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.52%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'LEO"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.353"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.04%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'-0.98%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.564"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-6.35%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'1.90%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2761"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.90%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04070"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.39%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'2.34%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004102"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.61%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001301"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.04%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D38").Value = "'0.02417"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'1.88%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05196"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.29%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007775"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-2.26%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1303"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.87%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007048"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'11.27%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.001952"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-1.55%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008806"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'5.69%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3330"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.29%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006958"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'6.01%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.06%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.03066"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'1,286.89%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.06%"
$ws.Range("E51").Style = "Normal"
